# Cluster analysis Fe-number Frost - update control point counts
# - B2: 2106 -> 2152
# - B3: 1608 -> 2106
# - Remove row 4 (A4=2, B4=544) entirely, shrinking the used range to A1:B3

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 2152
$ws.Range("B3").Value = 2106

# Delete the now-obsolete 4th row (bin "2" with count 544) so the sheet
# dimension shrinks from A1:B4 to A1:B3
$ws.Rows.Item(4).Delete()
